$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# batch3: gender for five subjects was mis-coded as "w"; fix to "f"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("batch3")
$ws3.Range("C1").Value = "f"
$ws3.Range("C3").Value = "f"
$ws3.Range("C4").Value = "f"
$ws3.Range("C5").Value = "f"
$ws3.Range("C9").Value = "f"
$ws3.Range("C1:C1048576").Select()

# ---------------------------------------------------------------------------
# batch4: fix the same gender typo, then add the last two recruited subjects
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("batch4")
$ws4.Range("C2").Value = "f"
$ws4.Range("C3").Value = "f"

$ws4.Range("A4").Value = 1706
$ws4.Range("B4").Value = 72
$ws4.Range("C4").Value = "f"
$ws4.Range("D4").Value = "PD+"
$ws4.Range("E4").Value = 10
$ws4.Range("F4").Value = 28
$ws4.Range("G4").Value = 2
$ws4.Range("H4").Value = 0.33333333333333331
$ws4.Range("H4").NumberFormat = "h:mm"
$ws4.Range("I4").Value = 0.5
$ws4.Range("I4").NumberFormat = "h:mm"
$ws4.Range("J4").Value = 0.40625
$ws4.Range("J4").NumberFormat = "h:mm"

$ws4.Range("A5").Value = 1722
$ws4.Range("B5").Value = 74
$ws4.Range("C5").Value = "f"
$ws4.Range("D5").Value = "PD-"
$ws4.Range("E5").Value = 7
$ws4.Range("F5").Value = 29
$ws4.Range("G5").Value = 1
$ws4.Range("H5").Value = 0.5
$ws4.Range("H5").NumberFormat = "h:mm"
$ws4.Range("I5").Value = 0.66666666666666663
$ws4.Range("I5").NumberFormat = "h:mm"
$ws4.Range("J5").Value = 0.67708333333333337
$ws4.Range("J5").NumberFormat = "h:mm"

$ws4.Range("H6").NumberFormat = "h:mm"
$ws4.Range("I6").NumberFormat = "h:mm"
$ws4.Range("J6").NumberFormat = "h:mm"
$ws4.Range("K6").NumberFormat = "h:mm"

$ws4.Range("C6").Select()

# ---------------------------------------------------------------------------
# all: consolidated sheet gets the same two new rows appended at the bottom
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("all")

$ws5.Range("A104").Value = 1706
$ws5.Range("B104").Value = 72
$ws5.Range("C104").Value = "f"
$ws5.Range("D104").Value = "PD+"
$ws5.Range("E104").Value = 10
$ws5.Range("F104").Value = 28
$ws5.Range("G104").Value = 2
$ws5.Range("H104").Value = 0.33333333333333331
$ws5.Range("H104").NumberFormat = "h:mm"
$ws5.Range("I104").Value = 0.5
$ws5.Range("I104").NumberFormat = "h:mm"
$ws5.Range("J104").Value = 0.40625
$ws5.Range("J104").NumberFormat = "h:mm"

$ws5.Range("A105").Value = 1722
$ws5.Range("B105").Value = 74
$ws5.Range("C105").Value = "f"
$ws5.Range("D105").Value = "PD-"
$ws5.Range("E105").Value = 7
$ws5.Range("F105").Value = 29
$ws5.Range("G105").Value = 1
$ws5.Range("H105").Value = 0.5
$ws5.Range("H105").NumberFormat = "h:mm"
$ws5.Range("I105").Value = 0.66666666666666663
$ws5.Range("I105").NumberFormat = "h:mm"
$ws5.Range("J105").Value = 0.67708333333333337
$ws5.Range("J105").NumberFormat = "h:mm"

$ws5.Range("R99").Select()
$ws5.Activate()
